# feat: add 2022-Q1 data
#
# 1) Create a new worksheet "2022-Q1" positioned after "2021-Q4" and before
#    "总计", with the same layout/style as the other per-quarter sheets.
# 2) Insert a matching summary row into the "总计" sheet, shifting the
#    existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet
# ---------------------------------------------------------------------
# Duplicate "2021-Q4" (same header set incl. "基金规模" + border/bold style)
# right after itself, then rename it and drop its extra data row so only
# a single fund remains, matching the target content.
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"
$newSheet.Rows("3:3").Delete()

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'004397"
$newSheet.Range("C2").Value = "长盛信息安全量化策略灵活配置混合"
$newSheet.Range("D2").Value = "'4.21"
$newSheet.Range("E2").Value = "'29.75"
$newSheet.Range("F2").Value = "'0.70"
$newSheet.Range("G2").Value = "'0.0295"
$newSheet.Range("H2").Value = 8

# ---------------------------------------------------------------------
# 2. Update "总计" sheet with the new 2022-Q1 summary row
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows("2:2").Insert()

# Row-insert copies the bold header style down from row 1 into the new
# row's B:D cells; restore the plain data-row styling (and the
# bordered/bold index-column style for A2) by copying formats from the
# row right below, which still holds the original row-2 formatting.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.03

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
